# Tripadvisor New Orleans shard 205 update:
#  1. Insert a new "State" column into hotel_info (between Hotel_Name and
#     City) and fill it with "Louisiana" for the existing hotel row.
#  2. Reorder the worksheet tabs so review_info comes before hotel_info.

$wb = $excel.ActiveWorkbook

# --- 1. Add the State column to hotel_info ---------------------------------
$wsHotel = $wb.Worksheets.Item("hotel_info")

# Hotel_Name is column B, City is column C -> insert a new column C so the
# new column sits between them, then push the header + data in.
$wsHotel.Columns.Item(3).Insert()
$wsHotel.Cells.Item(1, 3).Value = "State"
$wsHotel.Cells.Item(2, 3).Value = "Louisiana"

# --- 2. Reorder sheets: review_info, then hotel_info ------------------------
$wsReview = $wb.Worksheets.Item("review_info")
$wsReview.Move($wsHotel)
